# Commit: "completed citations for websites"
#
# Adds a final citation paragraph for the Creative Commons license page,
# right after the existing "instructables" citation block, matching the
# style of the other link citations already in the document:
#   <hyperlink to the url>, then a literal trailing space, as its own
#   paragraph, inserted just before the two trailing blank paragraphs.

$d = $word.ActiveDocument

# Locate the paragraph that contains the "instructables" citation link;
# the new citation goes right after the first (empty) paragraph that
# follows it.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*instructables.com*") {
        $anchorIndex = $i + 1
        break
    }
}
if ($anchorIndex -eq -1) {
    $anchorIndex = 8
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$urlText = "https://creativecommons.org/licenses/by-sa/3.0/"

# Type the URL text into the new (empty) paragraph, ahead of its
# paragraph mark, followed by a literal trailing space.
$startRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$startRange.InsertAfter($urlText)

$spaceRange = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$spaceRange.InsertAfter(" ")

# Turn just the URL text (not the trailing space) into a live hyperlink,
# same as the other citations in the document.
$urlRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $urlText.Length)
$d.Hyperlinks.Add($urlRange, $urlText)

Write-Output "Inserted citation paragraph: [$($newPara.Range.Text)]"
